# The deck's design is switched from the "Integral" theme to the stock
# built-in "Office Theme" palette. The font scheme and format scheme
# (fills/lines/effects) of "Integral" and "Office Theme" are identical in
# this deck, so the only thing that actually changes is the twelve
# scheme colours (Integral's greens/golds -> Office's blues/greys).

$p = $ppt.ActivePresentation

$theme  = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# ThemeColorScheme index order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink. .RGB is an OLE colour int (0xBBGGRR) - values below are the
# standard Office theme colours 000000/FFFFFF/44546A/E7E6E6/5B9BD5/
# ED7D31/A5A5A5/FFC000/4472C4/70AD47/0563C1/954F72.
$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
